$wb = $excel.ActiveWorkbook

# ALC!row2
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 369.57144
$ws.Cells.Item(2, 9).Value = 200
$ws.Cells.Item(2, 11).Value = 200
$ws.Cells.Item(2, 13).Value = -87

# ALC!row6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1054.8182
$ws.Cells.Item(6, 9).Value = 325.375
$ws.Cells.Item(6, 11).Value = 976.125
$ws.Cells.Item(6, 13).Value = -864.125

# ALC!row8
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 9.6
$ws.Cells.Item(8, 9).Value = 9.6
$ws.Cells.Item(8, 11).Value = 28.8
$ws.Cells.Item(8, 13).Value = 110.2

# ALC!row11
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 25.5
$ws.Cells.Item(11, 9).Value = 25.5
$ws.Cells.Item(11, 11).Value = 25.5
$ws.Cells.Item(11, 13).Value = 114.5

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2693.3872
$ws.Cells.Item(17, 10).Value = 2796.111
$ws.Cells.Item(17, 12).Value = 8388.332999999999
$ws.Cells.Item(17, 14).Value = -8724.332999999999

# ALC!row53
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 806.6667
$ws.Cells.Item(53, 9).Value = 495.5
$ws.Cells.Item(53, 10).Value = 1195.625
$ws.Cells.Item(53, 11).Value = 495.5
$ws.Cells.Item(53, 12).Value = 1195.625
$ws.Cells.Item(53, 13).Value = 141.5
$ws.Cells.Item(53, 14).Value = -2469.625

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 9499.286
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 9499.286
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 28497.858
$ws.Cells.Item(70, 14).Value = -29037.858
$ws.Cells.Item(70, 13).ClearContents()

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 9499.286
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 9499.286
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 28497.858
$ws.Cells.Item(73, 14).Value = -30369.858
$ws.Cells.Item(73, 13).ClearContents()

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 1000
$ws.Cells.Item(76, 9).Value = 1000
$ws.Cells.Item(76, 11).Value = 1000
$ws.Cells.Item(76, 13).Value = -685

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 1000
$ws.Cells.Item(79, 9).Value = 1000
$ws.Cells.Item(79, 11).Value = 1000
$ws.Cells.Item(79, 13).Value = 92

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 700
$ws.Cells.Item(86, 10).Value = 700
$ws.Cells.Item(86, 12).Value = 700
$ws.Cells.Item(86, 14).Value = -2946

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 700
$ws.Cells.Item(89, 10).Value = 700
$ws.Cells.Item(89, 12).Value = 3500
$ws.Cells.Item(89, 14).Value = -14732

# ALC!row92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1408.5454
$ws.Cells.Item(92, 9).Value = 143
$ws.Cells.Item(92, 11).Value = 143
$ws.Cells.Item(92, 13).Value = 1105

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 180.2
$ws.Cells.Item(107, 9).Value = 180.2
$ws.Cells.Item(107, 11).Value = 180.2
$ws.Cells.Item(107, 13).Value = 1739.8

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3807.2856
$ws.Cells.Item(122, 9).Value = 3858.5
$ws.Cells.Item(122, 11).Value = 11575.5
$ws.Cells.Item(122, 13).Value = -9125.5

# BSM!row5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 241.5
$ws.Cells.Item(5, 9).Value = 241.5
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 241.5
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -128.5
$ws.Cells.Item(5, 14).ClearContents()

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1724.3334
$ws.Cells.Item(105, 9).Value = 1667
$ws.Cells.Item(105, 11).Value = 1667
$ws.Cells.Item(105, 13).Value = 80

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2791.7646
$ws.Cells.Item(107, 9).Value = 843.0769
$ws.Cells.Item(107, 11).Value = 843.0769
$ws.Cells.Item(107, 13).Value = 1076.9231

# CRP!row22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1312.4445
$ws.Cells.Item(22, 9).Value = 787.75
$ws.Cells.Item(22, 10).Value = 2361.8333
$ws.Cells.Item(22, 11).Value = 787.75
$ws.Cells.Item(22, 12).Value = 2361.8333
$ws.Cells.Item(22, 13).Value = -437.75
$ws.Cells.Item(22, 14).Value = -3061.8333

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).ClearContents()

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).ClearContents()

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 969.4667
$ws.Cells.Item(122, 10).Value = 1036.75
$ws.Cells.Item(122, 12).Value = 3110.25
$ws.Cells.Item(122, 14).Value = -8010.25

# CUL!row6
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 125
$ws.Cells.Item(6, 10).Value = 166.66667
$ws.Cells.Item(6, 12).Value = 500.00001
$ws.Cells.Item(6, 14).Value = -726.00001

# CUL!row7
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 46.25
$ws.Cells.Item(7, 9).Value = 38.57143
$ws.Cells.Item(7, 11).Value = 115.71429
$ws.Cells.Item(7, 13).Value = -3.714290000000005

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 3123
$ws.Cells.Item(34, 10).Value = 3647.6
$ws.Cells.Item(34, 12).Value = 10942.8
$ws.Cells.Item(34, 14).Value = -11110.8

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 6750
$ws.Cells.Item(39, 10).Value = 7940
$ws.Cells.Item(39, 12).Value = 23820
$ws.Cells.Item(39, 14).Value = -24408

# CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 4585.5713
$ws.Cells.Item(55, 10).Value = 5799.8
$ws.Cells.Item(55, 12).Value = 17399.4
$ws.Cells.Item(55, 14).Value = -17753.4

# CUL!row121
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1495.9166
$ws.Cells.Item(121, 10).Value = 2428.4285
$ws.Cells.Item(121, 12).Value = 7285.2855
$ws.Cells.Item(121, 14).Value = -9905.2855

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2186.25
$ws.Cells.Item(131, 9).Value = 1922.5
$ws.Cells.Item(131, 10).Value = 2450
$ws.Cells.Item(131, 11).Value = 5767.5
$ws.Cells.Item(131, 12).Value = 7350
$ws.Cells.Item(131, 13).Value = -727.5
$ws.Cells.Item(131, 14).Value = -17430

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7002.2
$ws.Cells.Item(70, 9).Value = 7603
$ws.Cells.Item(70, 10).Value = 6401.4
$ws.Cells.Item(70, 11).Value = 7603
$ws.Cells.Item(70, 12).Value = 6401.4
$ws.Cells.Item(70, 13).Value = -7333
$ws.Cells.Item(70, 14).Value = -6941.4

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 7002.2
$ws.Cells.Item(73, 9).Value = 7603
$ws.Cells.Item(73, 10).Value = 6401.4
$ws.Cells.Item(73, 11).Value = 7603
$ws.Cells.Item(73, 12).Value = 6401.4
$ws.Cells.Item(73, 13).Value = -6667
$ws.Cells.Item(73, 14).Value = -8273.4

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 27699.12
$ws.Cells.Item(132, 9).Value = 35715.805
$ws.Cells.Item(132, 10).Value = 5106.636
$ws.Cells.Item(132, 11).Value = 107147.415
$ws.Cells.Item(132, 12).Value = 15319.908
$ws.Cells.Item(132, 13).Value = -104617.415
$ws.Cells.Item(132, 14).Value = -20379.908

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2380.5
$ws.Cells.Item(55, 9).Value = 2790
$ws.Cells.Item(55, 10).Value = 1971
$ws.Cells.Item(55, 11).Value = 2790
$ws.Cells.Item(55, 12).Value = 1971
$ws.Cells.Item(55, 13).Value = -2617
$ws.Cells.Item(55, 14).Value = -2317

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 5070.125
$ws.Cells.Item(100, 9).Value = 2531.3333
$ws.Cells.Item(100, 10).Value = 6593.4
$ws.Cells.Item(100, 11).Value = 2531.3333
$ws.Cells.Item(100, 12).Value = 6593.4
$ws.Cells.Item(100, 13).Value = -1990.3333

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3283.8333
$ws.Cells.Item(132, 9).Value = 2609.6
$ws.Cells.Item(132, 10).Value = 4126.625
$ws.Cells.Item(132, 11).Value = 7828.799999999999
$ws.Cells.Item(132, 12).Value = 12379.875
$ws.Cells.Item(132, 13).Value = -5298.799999999999
$ws.Cells.Item(132, 14).Value = -17439.875
